$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "7.2.1" labels to "7.2.1.1" in the header row (A1:C1)
$ws.Range("A1").Value = " 7.2.1.1 Энергия керектөөлөрүнүн жалпы көлөмүндөгү энергиянын жаңыланган булактарынын  үлүшү"
$ws.Range("B1").Value = " 7.2.1.1 Доля возобновляемых источников энергии в общем объеме энергопотребления"
$ws.Range("C1").Value = "7.2.1.1 Renewable energy share in the total energy consumption"

# 2. Update the selected cell in the sheet view from P9 to P7
$ws.Range("P7").Select()

# 3. Fill in the previously empty Q5 cell
$ws.Range("Q5").Value = 36.700000000000003

# 4. Update P6 and Q6 values
$ws.Range("P6").Value = 13859.2
$ws.Range("Q6").Value = 13979.2
